# Generate Report for Handoff
# Status moved from "In Translation" to "Ready for handoff" and the
# handoff/generate timestamps were refreshed. Updating the displayed
# status text on all three sheets causes the Status/zh-cn/de-de columns
# to re-autofit to the new (wider) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps ---
$wsOverview.Range("G2").Value = "2016-09-01 13:06:25"
$wsZhCn.Range("H2").Value = "2016-09-01 13:06:21"
$wsDeDe.Range("H2").Value = "2016-09-01 13:06:25"

# --- Columns widen to fit the new, longer status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
